$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; existing B,C shift to C,D.
$null = $ws.Columns("B:B").Insert()

# New column B gets a "StatQuery" header/value pair (a second Neo4j query
# column added alongside the existing query/dbExcel/WebExcel ones), for the
# diagnosis/gender/race/ethnicity test case stats.
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.race IN ['UNKNOWN'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Match column B's width to column A's width (both ~75.8 chars wide).
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Selection moves to single cell A4 instead of the old B2:B4 / B4 active cell.
$null = $ws.Range("A4").Select()

Write-Host "done"
